# Daily attendance processing - 2025-10-24 18:27:49
# Reorders the "Recorded By" (column G) names in the Session Analysis
# Results sheet: moves "System"/"system" from the front of the list to
# after the human recorder(s).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows where column G currently reads "system, backup@backdoor.com, System"
# and should become "backup@backdoor.com, System, system"
$rowsSystemFirst = @(2, 29, 56)
foreach ($r in $rowsSystemFirst) {
    $ws.Cells.Item($r, 7).Value = "backup@backdoor.com, System, system"
}

# Rows where column G currently reads "dnasr281@gmail.com, System"
# and should become "System, dnasr281@gmail.com"
$rowsDnasrFirst = @(3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 30, 33, 37, 39, 40, 41, 42, 45, 46, 47, 48, 49, 57, 60, 64, 66, 67, 68, 69, 72, 73, 74, 75, 76, 86, 87, 88, 89, 93, 95, 102, 112, 113, 114, 115, 119, 121, 128, 138, 139, 140, 141, 145, 147, 154)
foreach ($r in $rowsDnasrFirst) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}
